# Apply cryptos list value updates (price & volume columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $style = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $style
}

Set-TextValue $ws.Range("D2") "70.773.95"
Set-TextValue $ws.Range("E2") "  +2.68%  "
Set-TextValue $ws.Range("D3") "3.796.83"
Set-TextValue $ws.Range("E3") "  +0.82%  "
Set-TextValue $ws.Range("E4") "  -0.16%  "
Set-TextValue $ws.Range("D5") "700.61"
Set-TextValue $ws.Range("E5") "  +8.90%  "
Set-TextValue $ws.Range("D6") "172.59"
Set-TextValue $ws.Range("E6") "  +4.32%  "
Set-TextValue $ws.Range("D7") "3.796.59"
Set-TextValue $ws.Range("E7") "  +0.88%  "
Set-TextValue $ws.Range("E8") "  -0.02%  "
Set-TextValue $ws.Range("D9") "0.529"
Set-TextValue $ws.Range("E9") "  +1.23%  "
Set-TextValue $ws.Range("E10") "  +2.77%  "
Set-TextValue $ws.Range("E11") "  +6.48%  "
Set-TextValue $ws.Range("E12") "  +0.94%  "
Set-TextValue $ws.Range("D13") "0.0000258"
Set-TextValue $ws.Range("E13") "  +8.22%  "
Set-TextValue $ws.Range("D14") "36.42"
Set-TextValue $ws.Range("E14") "  +4.46%  "
Set-TextValue $ws.Range("D15") "4.435.13"
Set-TextValue $ws.Range("E15") "  +0.83%  "
Set-TextValue $ws.Range("D16") "3.795.44"
Set-TextValue $ws.Range("E16") "  +0.72%  "
Set-TextValue $ws.Range("D17") "70.775.59"
Set-TextValue $ws.Range("E17") "  +2.69%  "
Set-TextValue $ws.Range("D18") "17.90"
Set-TextValue $ws.Range("E18") "  +1.40%  "
Set-TextValue $ws.Range("D19") "7.21"
Set-TextValue $ws.Range("E19") "  +3.09%  "
Set-TextValue $ws.Range("D21") "11.10"
Set-TextValue $ws.Range("E21") "  +16.02%  "
Set-TextValue $ws.Range("D22") "482.48"
Set-TextValue $ws.Range("E22") "  +2.11%  "
Set-TextValue $ws.Range("E23") "  +1.40%  "
Set-TextValue $ws.Range("D24") "84.24"
Set-TextValue $ws.Range("E24") "  +3.04%  "
Set-TextValue $ws.Range("E25") "  +0.76%  "
Set-TextValue $ws.Range("D26") "12.46"
Set-TextValue $ws.Range("E26") "  +2.56%  "
Set-TextValue $ws.Range("E27") "  +3.89%  "
Set-TextValue $ws.Range("D28") "10.46"
Set-TextValue $ws.Range("E28") "  +4.04%  "
Set-TextValue $ws.Range("D29") "3.946.41"
Set-TextValue $ws.Range("E29") "  +0.80%  "
Set-TextValue $ws.Range("E30") "  -0.12%  "
Set-TextValue $ws.Range("E31") "  +16.49%  "
Set-TextValue $ws.Range("E32") "  +5.62%  "
Set-TextValue $ws.Range("D33") "2.29"
Set-TextValue $ws.Range("E33") "  +0.83%  "
Set-TextValue $ws.Range("D34") "0.184"
Set-TextValue $ws.Range("E34") "  +6.60%  "
Set-TextValue $ws.Range("E35") "  +3.46%  "
Set-TextValue $ws.Range("E36") "  +3.89%  "
Set-TextValue $ws.Range("D37") "0.999"
Set-TextValue $ws.Range("E37") "  +0.02%  "
Set-TextValue $ws.Range("E38") "  +2.09%  "
Set-TextValue $ws.Range("E39") "  +6.42%  "
Set-TextValue $ws.Range("D40") "6.05"
Set-TextValue $ws.Range("E40") "  +4.89%  "
Set-TextValue $ws.Range("D41") "2.20"
Set-TextValue $ws.Range("E41") "  +11.56%  "
Set-TextValue $ws.Range("E42") "  +1.91%  "
Set-TextValue $ws.Range("D43") "0.000327"
Set-TextValue $ws.Range("E43") "  +23.06%  "
Set-TextValue $ws.Range("D44") "1.00"
Set-TextValue $ws.Range("E44") "  -0.08%  "
Set-TextValue $ws.Range("E45") "  +0.01%  "
Set-TextValue $ws.Range("D46") "45.73"
Set-TextValue $ws.Range("E46") "  +1.38%  "
Set-TextValue $ws.Range("D47") "162.12"
Set-TextValue $ws.Range("E47") "  +4.27%  "
Set-TextValue $ws.Range("D48") "48.78"
Set-TextValue $ws.Range("E48") "  +1.96%  "
Set-TextValue $ws.Range("D49") "0.302"
Set-TextValue $ws.Range("E49") "  +2.39%  "
Set-TextValue $ws.Range("E50") "  -1.24%  "
Set-TextValue $ws.Range("D51") "8.59"
Set-TextValue $ws.Range("E51") "  +2.61%  "
